$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add detail on gas limit (row 17, col D) and block interval (row 18, col D)
$ws.Range("D17").Value = "20s block interval"
$ws.Range("D18").Value = "100millions Gas Limit"

# Update the active selection to match D18 (as in diff)
$ws.Range("D18").Select()
